# Normalize the "Recorded By" (column G) values on the Session Analysis
# sheet: each cell holds a comma-separated list of recorder names/emails
# that needs to be sorted using an ordinal (case-sensitive, ASCII-code)
# string sort, e.g.:
#   "system, System, backup@backdoor.com" -> "System, backup@backdoor.com, system"
#   "dnasr281@gmail.com, System"          -> "System, dnasr281@gmail.com"
#   "dnasr281@gmail.com, admin@admin.com" -> "admin@admin.com, dnasr281@gmail.com"
#
# NOTE: this runtime's PS-like interpreter does not give function calls
# their own variable scope (writes inside a function leak into the
# caller's scope), so every helper below uses uniquely-prefixed variable
# names to avoid clobbering the loop counters of its callers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Compare-Ordinal($cmpX, $cmpY) {
    $cmpLenX = $cmpX.Length
    $cmpLenY = $cmpY.Length
    $cmpMinLen = $cmpLenX
    if ($cmpLenY -lt $cmpMinLen) { $cmpMinLen = $cmpLenY }
    $cmpIdx = 0
    $cmpResult = 0
    while ($cmpIdx -lt $cmpMinLen) {
        $cmpCx = [int][char]$cmpX[$cmpIdx]
        $cmpCy = [int][char]$cmpY[$cmpIdx]
        if ($cmpCx -lt $cmpCy) { $cmpResult = -1; $cmpIdx = $cmpMinLen }
        elseif ($cmpCx -gt $cmpCy) { $cmpResult = 1; $cmpIdx = $cmpMinLen }
        else { $cmpIdx = $cmpIdx + 1 }
    }
    if ($cmpResult -eq 0) {
        if ($cmpLenX -lt $cmpLenY) { $cmpResult = -1 }
        elseif ($cmpLenX -gt $cmpLenY) { $cmpResult = 1 }
    }
    return $cmpResult
}

function Sort-Ordinal($srtArr) {
    $srtN = $srtArr.Count
    $srtList = @()
    foreach ($srtX in $srtArr) { $srtList += $srtX }
    $srtI = 0
    while ($srtI -lt $srtN) {
        $srtJ = 0
        while ($srtJ -lt ($srtN - $srtI - 1)) {
            $srtCmp = Compare-Ordinal $srtList[$srtJ] $srtList[$srtJ+1]
            if ($srtCmp -gt 0) {
                $srtTmp = $srtList[$srtJ]
                $srtList[$srtJ] = $srtList[$srtJ+1]
                $srtList[$srtJ+1] = $srtTmp
            }
            $srtJ = $srtJ + 1
        }
        $srtI = $srtI + 1
    }
    return $srtList
}

function Get-SortedRecordedBy($rbyText) {
    $rbyParts = $rbyText.Split(",")
    $rbyTrimmed = @()
    foreach ($rbyP in $rbyParts) {
        $rbyTrimmed += $rbyP.Trim()
    }
    $rbySorted = Sort-Ordinal $rbyTrimmed
    $rbyJoined = [string]::Join(", ", $rbySorted)
    return $rbyJoined
}

$colRecordedBy = 7
$firstDataRow = 2
$lastDataRow = $ws.UsedRange.Rows.Count

$row = $firstDataRow
while ($row -le $lastDataRow) {
    $cell = $ws.Cells.Item($row, $colRecordedBy)
    $origValue = $cell.Text
    if ($origValue -ne "") {
        $newValue = Get-SortedRecordedBy $origValue
        if ($newValue -ne $origValue) {
            $cell.Value = $newValue
        }
    }
    $row = $row + 1
}
